$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF), styled like the other
# header cells in row 1 (same formatting as H1 - bold, bordered, centered).
$xlPasteFormats = -4122

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J1").Value = "IF"

$excel.CutCopyMode = $false

# Fill data rows 2-30: column I is always 1, column J mirrors column H.
for ($row = 2; $row -le 30; $row++) {
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $ws.Cells.Item($row, 8).Value2
}
